$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.816.99"
$ws.Range("E2").Value = "  -4.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.674.19"
$ws.Range("E3").Value = "  -5.25%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.73"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.97"
$ws.Range("E6").Value = "  +6.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.665.40"
$ws.Range("E7").Value = "  -5.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("E8").Value = "  -7.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.710"
$ws.Range("E10").Value = "  -6.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  -8.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.23"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("E13").Value = "  -10.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.30"
$ws.Range("E14").Value = "  -10.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.270.67"
$ws.Range("E15").Value = "  -5.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.689.24"
$ws.Range("E16").Value = "  -5.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.26"
$ws.Range("E17").Value = "  -8.37%  "
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  -8.55%  "
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.11"
$ws.Range("E20").Value = "  -7.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.636.09"
$ws.Range("E21").Value = "  -4.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "406.06"
$ws.Range("E22").Value = "  -7.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.50"
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.14"
$ws.Range("E24").Value = "  -7.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.99"
$ws.Range("E25").Value = "  -8.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.68"
$ws.Range("E26").Value = "  -8.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.83"
$ws.Range("E28").Value = "  -7.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.05"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.36"
$ws.Range("E30").Value = "  -9.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.49"
$ws.Range("E31").Value = "  -7.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.20"
$ws.Range("E32").Value = "  -13.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.40"
$ws.Range("E33").Value = "  -9.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.116"
$ws.Range("E34").Value = "  -7.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "43.14"
$ws.Range("E35").Value = "  -12.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.94"
$ws.Range("E36").Value = "  -8.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "584.74"
$ws.Range("E37").Value = "  -7.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0876"
$ws.Range("E38").Value = "  -11.86%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.395"
$ws.Range("E40").Value = "  -7.57%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").Value = "  -7.33%  "
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("E44").Value = "  -9.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0432"
$ws.Range("E45").Value = "  -8.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  -15.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.06"
$ws.Range("E47").Value = "  -11.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.69"
$ws.Range("E48").Value = "  -4.84%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.133"
$ws.Range("E49").Value = "  -7.64%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.738.14"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.14"
$ws.Range("E51").Value = "  -5.27%  "
